# Apply the betexplorer brazil_serie-b_2023 corrections:
#  - A bunch of rows had their match-detail columns (F:V) shifted to the
#    wrong row (home/away teams, odds, timestamps, url). Fix by swapping /
#    rotating the F:V payloads back onto the correct row (A:E - Indice,
#    pais, torneio, temporada, data_partida - are untouched / already
#    correct for every row).
#  - 3 brand-new match rows are appended at the bottom (283-285).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Get-RowValues($r) {
    $result = @{}
    foreach ($c in $cols) {
        $addr = $c + $r
        $result[$c] = $ws.Range($addr).Value2
    }
    return $result
}

function Set-RowValues($r, $vals) {
    foreach ($c in $cols) {
        $addr = $c + $r
        $ws.Range($addr).Value = $vals[$c]
    }
}

# ---- simple pairwise swaps (F:V payload moves wholesale between the two
#      rows; A:E stay put) --------------------------------------------------
$swapPairs = @(
    @(57,58),
    @(64,65),
    @(74,76),
    @(85,86),
    @(115,116),
    @(135,136),
    @(156,157),
    @(163,164),
    @(174,175),
    @(182,183),
    @(194,195),
    @(203,204),
    @(223,224),
    @(235,236)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $v1 = Get-RowValues $r1
    $v2 = Get-RowValues $r2
    Set-RowValues $r1 $v2
    Set-RowValues $r2 $v1
}

# ---- 4-way rotations -------------------------------------------------------
# new(176) = old(178); new(177) = old(176); new(178) = old(179); new(179) = old(177)
$a176 = Get-RowValues 176
$a177 = Get-RowValues 177
$a178 = Get-RowValues 178
$a179 = Get-RowValues 179
Set-RowValues 176 $a178
Set-RowValues 177 $a176
Set-RowValues 178 $a179
Set-RowValues 179 $a177

# new(205) = old(207); new(207) = old(208); new(208) = old(209); new(209) = old(205)
$a205 = Get-RowValues 205
$a207 = Get-RowValues 207
$a208 = Get-RowValues 208
$a209 = Get-RowValues 209
Set-RowValues 205 $a207
Set-RowValues 207 $a208
Set-RowValues 208 $a209
Set-RowValues 209 $a205

# ---- append 3 new rows at the bottom (283, 284, 285) -----------------------
# Copy the A/E cell formatting (styles 1 / 2) from an existing data row so the
# new rows render identically to the rest of the table.
$ws.Range("A282").Copy()
$ws.Range("A283:A285").PasteSpecial(-4122)
$ws.Range("E282").Copy()
$ws.Range("E283:E285").PasteSpecial(-4122)

function Set-NewRow($r, $indice, $data_partida, $home, $homeGols, $away, $awayGols, $hOpenOdds, $hOpenDt, $hCloseOdds, $hCloseDt, $drawOpenOdds, $drawOpenDt, $drawCloseOdds, $drawCloseDt, $aOpenOdds, $aOpenDt, $aCloseOdds, $aCloseDt, $url) {
    $ws.Range("A" + $r).Value = $indice
    $ws.Range("B" + $r).Value = "brazil"
    $ws.Range("C" + $r).Value = "serie-b"
    $ws.Range("D" + $r).Value = "2023"
    $ws.Range("E" + $r).Value = $data_partida
    $ws.Range("F" + $r).Value = $home
    $ws.Range("G" + $r).Value = $homeGols
    $ws.Range("H" + $r).Value = $away
    $ws.Range("I" + $r).Value = $awayGols
    $ws.Range("J" + $r).Value = $hOpenOdds
    $ws.Range("K" + $r).Value = $hOpenDt
    $ws.Range("L" + $r).Value = $hCloseOdds
    $ws.Range("M" + $r).Value = $hCloseDt
    $ws.Range("N" + $r).Value = $drawOpenOdds
    $ws.Range("O" + $r).Value = $drawOpenDt
    $ws.Range("P" + $r).Value = $drawCloseOdds
    $ws.Range("Q" + $r).Value = $drawCloseDt
    $ws.Range("R" + $r).Value = $aOpenOdds
    $ws.Range("S" + $r).Value = $aOpenDt
    $ws.Range("T" + $r).Value = $aCloseOdds
    $ws.Range("U" + $r).Value = $aCloseDt
    $ws.Range("V" + $r).Value = $url
}

Set-NewRow 283 282 45192 "Ponte Preta" 0 "Mirassol" 3 `
    2.94 "16/09/2023 19:43" 4.05 "22/09/2023 23:59" `
    2.75 "16/09/2023 19:43" 2.87 "22/09/2023 23:56" `
    2.94 "16/09/2023 19:43" 2.28 "22/09/2023 23:59" `
    "https://www.betexplorer.com/football/brazil/serie-b/ponte-preta-mirassol/bivEkNQ0/"

Set-NewRow 284 283 45192.10416666666 "Atletico GO" 3 "Criciuma" 1 `
    1.96 "16/09/2023 21:12" 1.91 "23/09/2023 01:54" `
    3.2 "16/09/2023 21:12" 3.26 "23/09/2023 01:50" `
    4.5 "16/09/2023 21:12" 4.88 "23/09/2023 01:54" `
    "https://www.betexplorer.com/football/brazil/serie-b/atletico-go-criciuma/tGTIlst7/"

Set-NewRow 285 284 45192.10416666666 "Ituano" 0 "Vitoria" 2 `
    2.61 "19/09/2023 10:42" 2.32 "23/09/2023 02:28" `
    2.88 "19/09/2023 10:42" 3.07 "23/09/2023 02:28" `
    3.2 "19/09/2023 10:42" 3.63 "23/09/2023 02:28" `
    "https://www.betexplorer.com/football/brazil/serie-b/ituano-vitoria/8ML4CpYD/"
